# Updated remaining queries for C3DC
# Fixes the LEFT JOIN conditions in every SQL query stored on Sheet1
# (columns B2:B7 and C2) so that joins use the real key columns
# (study_id / participant_id) instead of the generic "id" columns,
# and refreshes a couple of cosmetic view settings that changed
# alongside the content edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-Sql($text) {
    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    return $text
}

# Every cell on the sheet that holds one of the DuckDB-style queries
# whose join columns need correcting.
$cells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $rng.Value2 = Fix-Sql $rng.Value2
}

# Column C was widened (and is no longer auto "best fit") now that the
# updated query text needs a bit more room.
$ws.Columns("C").ColumnWidth = 73.1666666666667

# The active view had scrolled down one row and the active cell moved
# from C5 to C7.
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("C7").Select()
